# Update cryptos list data (price + 1h volume%) for the Fri Jul 28 07:32:05 UTC 2023 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.194.07'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.84%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.863.12'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.70%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7084'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.05'
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3087'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.71%  '
$ws.Range("E9").Value = '  -2.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.73'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08376'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.63%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.867.87'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.181'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7109'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.38'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.55%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.226.06'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.76%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.930'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '242.76'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.77%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007820'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.116.74'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.86%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.10'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.24%  '
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.873'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.24%  '
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("E25").Value = '  +1.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.27'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.925'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.46'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.498'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.407'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.80%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.310'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.224'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.40%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05133'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.8068'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +11.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.914'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.83%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.164'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.89%  '
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("E38").Value = '  -1.13%  '
$ws.Range("E39").Value = '  -0.99%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.168.61'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.82%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.182'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8898'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.94%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '72.90'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.60%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9999'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '102.15'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.013.95'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.30%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5183'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.57%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.771'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.08%  '
$ws.Range("E49").Value = '  +0.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.272'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.001'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.06%  '
